$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17 (columns A-T), per the natmi LR-pair
# re-run "following Dr Hou advice": the sending/target cluster set now
# includes an "M2" cluster alongside ECs/FAPs/sCs, so the Sema7a -> Plxnc1
# table grows from a 4x3 (12-row) to a 4x4 (16-row) cluster cross-product,
# and several underlying expression-summary numbers are recomputed.
$data = @(
    @("ECs","Sema7a","Plxnc1","ECs",3,1,19.76232966666667,59.28698900000001,0.6092018671069146,0.6092018671069145,3,1,16.75150633333334,50.254519,0.3469350574503514,0.3469350574503513,331.048790572588,2979.439115153291,0.2113534847635987,0.2113534847635987),
    @("ECs","Sema7a","Plxnc1","FAPs",3,1,19.76232966666667,59.28698900000001,0.6092018671069146,0.6092018671069145,3,1,4.930404,14.791212,0.1021120107622625,0.1021120107622625,97.43626923785202,876.926423140668,0.06220682761041169,0.06220682761041166),
    @("ECs","Sema7a","Plxnc1","M2",3,1,19.76232966666667,59.28698900000001,0.6092018671069146,0.6092018671069145,3,1,25.470675,76.412025,0.5275149540934355,0.5275149540934354,503.3598761825251,4530.238885642725,0.3213630949605393,0.3213630949605391),
    @("ECs","Sema7a","Plxnc1","sCs",3,1,19.76232966666667,59.28698900000001,0.6092018671069146,0.6092018671069145,3,1,1.131685666666667,3.395057,0.02343797769395061,0.02343797769395061,22.36474522370811,201.282707013373,0.01427845977236493,0.01427845977236493),
    @("FAPs","Sema7a","Plxnc1","ECs",3,1,9.791583333333334,29.37475,0.3018394566436632,0.3018394566436631,3,1,16.75150633333334,50.254519,0.3469350574503514,0.3469350574503513,164.0237702216945,1476.21393199525,0.1047186892314521,0.1047186892314521),
    @("FAPs","Sema7a","Plxnc1","FAPs",3,1,9.791583333333334,29.37475,0.3018394566436632,0.3018394566436631,3,1,4.930404,14.791212,0.1021120107622625,0.1021120107622625,48.276461633,434.488154697,0.03082143384527321,0.03082143384527319),
    @("FAPs","Sema7a","Plxnc1","M2",3,1,9.791583333333334,29.37475,0.3018394566436632,0.3018394566436631,3,1,25.470675,76.412025,0.5275149540934355,0.5275149540934354,249.39823681875,2244.58413136875,0.1592248271149695,0.1592248271149694),
    @("FAPs","Sema7a","Plxnc1","sCs",3,1,9.791583333333334,29.37475,0.3018394566436632,0.3018394566436631,3,1,1.131685666666667,3.395057,0.02343797769395061,0.02343797769395061,11.08099451230555,99.72895061074998,0.007074506451968351,0.007074506451968349),
    @("M2","Sema7a","Plxnc1","ECs",1,0.3333333333333333,0.05307799999999999,0.159234,0.001636204700948844,0.001636204700948844,3,1,16.75150633333334,50.254519,0.3469350574503514,0.3469350574503513,0.8891364531606667,8.002228078446,0.0005676567719242222,0.0005676567719242221),
    @("M2","Sema7a","Plxnc1","FAPs",1,0.3333333333333333,0.05307799999999999,0.159234,0.001636204700948844,0.001636204700948844,3,1,4.930404,14.791212,0.1021120107622625,0.1021120107622625,0.261695983512,2.355263851608,0.0001670761520325529,0.0001670761520325528),
    @("M2","Sema7a","Plxnc1","M2",1,0.3333333333333333,0.05307799999999999,0.159234,0.001636204700948844,0.001636204700948844,3,1,25.470675,76.412025,0.5275149540934355,0.5275149540934354,1.35193248765,12.16739238885,0.0008631224477084929,0.0008631224477084926),
    @("M2","Sema7a","Plxnc1","sCs",1,0.3333333333333333,0.05307799999999999,0.159234,0.001636204700948844,0.001636204700948844,3,1,1.131685666666667,3.395057,0.02343797769395061,0.02343797769395061,0.06006761181533332,0.5406085063379998,[double]"3.834932928357614E-05",[double]"3.834932928357614E-05"),
    @("sCs","Sema7a","Plxnc1","ECs",3,1,2.832715333333333,8.498146,0.08732247154847343,0.08732247154847342,3,1,16.75150633333334,50.254519,0.3469350574503514,0.3469350574503513,47.45224884686378,427.070239621774,0.0302952266833763,0.03029522668337629),
    @("sCs","Sema7a","Plxnc1","FAPs",3,1,2.832715333333333,8.498146,0.08732247154847343,0.08732247154847342,3,1,4.930404,14.791212,0.1021120107622625,0.1021120107622625,13.966431010328,125.697879092952,0.008916673154545082,0.008916673154545078),
    @("sCs","Sema7a","Plxnc1","M2",3,1,2.832715333333333,8.498146,0.08732247154847343,0.08732247154847342,3,1,25.470675,76.412025,0.5275149540934355,0.5275149540934354,72.15117162285,649.36054460565,0.04606390957021829,0.04606390957021828),
    @("sCs","Sema7a","Plxnc1","sCs",3,1,2.832715333333333,8.498146,0.08732247154847343,0.08732247154847342,3,1,1.131685666666667,3.395057,0.02343797769395061,0.02343797769395061,3.205743340480222,28.851690064322,0.002046662140333758,0.002046662140333757)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}

Write-Host "Wrote $($data.Length) data rows"
